# Insert two new weekly price rows ("Hortaliza, Terminal Hortofrutícola
# Agro Chillán - Sandia") above the existing last record, pushing the
# old row 160 down to row 162.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 160/161; everything from the old row 160
# onward (including the old row 160 itself) shifts down by two rows.
$ws.Range("A160:R161").EntireRow.Insert()

# New row 160: Sandia, Extra quality, Región de O'Higgins origin.
$ws.Cells.Item(160, 1).Value = 7
$ws.Cells.Item(160, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(160, 3).Value = "Ñuble"
$ws.Cells.Item(160, 4).Value = 44595
$ws.Cells.Item(160, 5).Value = 16
$ws.Cells.Item(160, 6).Value = 100112028
$ws.Cells.Item(160, 7).Value = "Sandia"
$ws.Cells.Item(160, 8).Value = "Sin especificar"
$ws.Cells.Item(160, 9).Value = "Extra"
$ws.Cells.Item(160, 10).Value = 300
$ws.Cells.Item(160, 11).Value = 2500
$ws.Cells.Item(160, 12).Value = 2500
$ws.Cells.Item(160, 13).Value = 2500
$ws.Cells.Item(160, 14).Value = "$/unidad"
$ws.Cells.Item(160, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(160, 16).Value = 2500
$ws.Cells.Item(160, 17).Value = 1
$ws.Cells.Item(160, 18).Value = "Hortaliza"

# New row 161: Sandia, Primera quality, Región de O'Higgins origin.
$ws.Cells.Item(161, 1).Value = 7
$ws.Cells.Item(161, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(161, 3).Value = "Ñuble"
$ws.Cells.Item(161, 4).Value = 44595
$ws.Cells.Item(161, 5).Value = 16
$ws.Cells.Item(161, 6).Value = 100112028
$ws.Cells.Item(161, 7).Value = "Sandia"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 400
$ws.Cells.Item(161, 11).Value = 1800
$ws.Cells.Item(161, 12).Value = 2000
$ws.Cells.Item(161, 13).Value = 1900
$ws.Cells.Item(161, 14).Value = "$/unidad"
$ws.Cells.Item(161, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(161, 16).Value = 1900
$ws.Cells.Item(161, 17).Value = 1
$ws.Cells.Item(161, 18).Value = "Hortaliza"
